# Fixed the assignment 4
# Remove the "Modular components" list item paragraph (a ListParagraph /
# numbered bullet) that followed "Feature descriptions ... in Word .docx
# format" and preceded "Grading".

$d = $word.ActiveDocument

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Modular components") {
        $p.Range.Delete()
        break
    }
}
